$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.084.71'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.809.01'
$ws.Range('E3').Value = '  +1.88%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '600.74'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '164.15'
$ws.Range('E6').Value = '  -2.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.810.28'
$ws.Range('E7').Value = '  +1.93%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -0.38%  '
$ws.Range('E10').Value = '  +1.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.33'
$ws.Range('E11').Value = '  -1.60%  '
$ws.Range('E12').Value = '  +0.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '36.99'
$ws.Range('E13').Value = '  -2.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000244'
$ws.Range('E14').Value = '  -0.64%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.445.83'
$ws.Range('E15').Value = '  +1.97%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.839.76'
$ws.Range('E16').Value = '  +2.82%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.253.32'
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.50'
$ws.Range('E18').Value = '  +2.98%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.52'
$ws.Range('E19').Value = '  +7.32%  '
$ws.Range('B20').Value = 'TRON'
$ws.Range('C20').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.114'
$ws.Range('E20').Value = '  +0.23%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.19'
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '486.65'
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('E23').Value = '  -0.62%  '
$ws.Range('E24').Value = '  +4.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.45'
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('E26').Value = '  -2.24%  '
$ws.Range('E27').Value = '  -1.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.05'
$ws.Range('E28').Value = '  -0.62%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('E30').Value = '  -0.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.99'
$ws.Range('E31').Value = '  -0.36%  '
$ws.Range('E32').Value = '  -4.27%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.963.16'
$ws.Range('E33').Value = '  +2.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '31.82'
$ws.Range('E34').Value = '  +1.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.752.19'
$ws.Range('E35').Value = '  +2.20%  '
$ws.Range('E36').Value = '  -1.67%  '
$ws.Range('E37').Value = '  +0.76%  '
$ws.Range('E38').Value = '  +4.45%  '
$ws.Range('E39').Value = '  +0.30%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('E41').Value = '  +2.80%  '
$ws.Range('E42').Value = '  -1.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '434.85'
$ws.Range('E43').Value = '  +0.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '48.54'
$ws.Range('E44').Value = '  +0.24%  '
$ws.Range('E45').Value = '  +0.32%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.36'
$ws.Range('E47').Value = '  -1.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '142.89'
$ws.Range('E48').Value = '  +0.85%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.824.13'
$ws.Range('E49').Value = '  +1.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0355'
$ws.Range('E50').Value = '  +0.37%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '39.21'
$ws.Range('E51').Value = '  -2.03%  '
